$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 376, pushing existing rows 376-403 down to 377-404.
$ws.Rows.Item(376).Insert()

# Populate the newly inserted row 376 with the new weekly data entry.
$ws.Cells.Item(376, 1).Value = 10
$ws.Cells.Item(376, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(376, 3).Value = "La Araucanía"
$ws.Cells.Item(376, 4).Value = 44578
$ws.Cells.Item(376, 5).Value = 9
$ws.Cells.Item(376, 6).Value = "Fruta"
$ws.Cells.Item(376, 7).Value = 100103
$ws.Cells.Item(376, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(376, 9).Value = 100103006
$ws.Cells.Item(376, 10).Value = "Nectarín"
$ws.Cells.Item(376, 11).Value = "Super Queen"
$ws.Cells.Item(376, 12).Value = "Especial"
$ws.Cells.Item(376, 13).Value = 215
$ws.Cells.Item(376, 14).Value = 23000
$ws.Cells.Item(376, 15).Value = 23000
$ws.Cells.Item(376, 16).Value = 23000
$ws.Cells.Item(376, 17).Value = "`$/caja 20 kilos empedrada"
$ws.Cells.Item(376, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(376, 19).Value = 1150
$ws.Cells.Item(376, 20).Value = 20
